$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C rows 2-45 hold the "Förändrad" (last changed) date, stored as a
# serial date number. The whole column of values was bumped by one day
# (45177 -> 45178, i.e. 2023-09-08 -> 2023-09-09).
$ws.Range("C2:C45").Value = 45178
